# Add an "Age" column between the participant-id column (A) and the
# Betas column (B), and give both B1/C1 proper header labels ("Age" /
# "Betas"). The existing Betas values (currently in column B) move to
# column C; column B gets the new Age values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Participant ages, in row order (rows 2-32), matching participant ids
# already present in column A.
$ages = @(74,75,67,70,80,70,67,73,74,72,76,71,70,69,68,73,73,70,66,67,86,73,66,69,75,72,75,71,61,69,70)

$lastRow = 32

# 1) Give the new header cell (C1) the same formatting (bold, bordered,
#    centered) as the existing header cell (B1) before we touch any
#    values.
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 2) Move the existing Betas values from column B down into column C,
#    working from the bottom up so we never overwrite a value we still
#    need to read.
for ($r = $lastRow; $r -ge 2; $r--) {
    $ws.Range("C$r").Value2 = $ws.Range("B$r").Value2
}

# 3) Fill column B with the new Age values.
for ($i = 0; $i -lt $ages.Length; $i++) {
    $r = $i + 2
    $ws.Range("B$r").Value2 = $ages[$i]
}

# 4) Set the header labels.
$ws.Range("B1").Value2 = "Age"
$ws.Range("C1").Value2 = "Betas"
